$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.772.95"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "2.213.70"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.23"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.49"
$ws.Range("E6").Value = "  +12.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.42"
$ws.Range("E10").Value = "  +10.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.58"
$ws.Range("E12").Value = "  +9.06%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "2.550.70"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.60"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "2.261.05"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "43.750.00"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.98"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.08"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("E22").Value = "  +8.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.07"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.02"
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.83"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.51"
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.49"
$ws.Range("E29").Value = "  -7.31%  "
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0892"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.53"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0360"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("E39").Value = "  +16.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.33"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.60"
$ws.Range("E41").Value = "  +7.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.09"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  +4.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.16"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0984"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.57"
$ws.Range("E49").Value = "  +9.35%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.447"
$ws.Range("E51").Value = "  -4.19%  "
